# Update crypto price/volume table to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.652.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.422.12'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.92'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.426.33'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.013.86'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000189'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.663.35'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.401.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.22'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.01'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.20%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.21'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.38%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.48'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.43%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.26'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.11'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.60'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.94'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.95'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.30%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0758'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.902.76'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.58'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.64'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.10'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0316'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.773'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.81'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '320.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.28%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.73%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.75%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.108'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.24%  '
